$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-20: column B gets the quoted .mat filename (the text itself
# contains literal single-quote characters at both ends), column C gets
# updated numeric values.
$data = @(
    @{Row=2;  Name="''Bacteroides_cellulosilyticus_DSM_14838.mat'";      Val=0}
    @{Row=3;  Name="''Bacteroides_coprocola_M16_DSM_17136.mat'";         Val=0}
    @{Row=4;  Name="''Bacteroides_coprophilus_DSM_18228.mat'";           Val=0}
    @{Row=5;  Name="''Bacteroides_fluxus_YIT_12057.mat'";                Val=0}
    @{Row=6;  Name="''Bacteroides_oleiciplenus_YIT_12058.mat'";          Val=0}
    @{Row=7;  Name="''Bacteroides_ovatus_ATCC_8483.mat'";                Val=0}
    @{Row=8;  Name="''Bacteroides_plebeius_M12_DSM_17135.mat'";          Val=0}
    @{Row=9;  Name="''Bacteroides_salyersiae_WAL_10018.mat'";            Val=0}
    @{Row=10; Name="''Bacteroides_stercoris_ATCC_43183.mat'";            Val=0}
    @{Row=11; Name="''Bacteroides_thetaiotaomicron_VPI_5482.mat'";       Val=0.016}
    @{Row=12; Name="''Bacteroides_uniformis_ATCC_8492.mat'";             Val=0}
    @{Row=13; Name="''Bacteroides_vulgatus_ATCC_8482.mat'";              Val=0.011}
    @{Row=14; Name="''Bifidobacterium_animalis_lactis_AD011.mat'";       Val=0}
    @{Row=15; Name="''Enterococcus_faecalis_OG1RF_ATCC_47077.mat'";      Val=0}
    @{Row=16; Name="''Flavonifractor_plautii_ATCC_29863.mat'";           Val=0}
    @{Row=17; Name="''Lactobacillus_plantarum_JDM1.mat'";                Val=0.014}
    @{Row=18; Name="''Odoribacter_laneus_YIT_12061.mat'";                Val=0.144}
    @{Row=19; Name="''Parabacteroides_distasonis_ATCC_8503.mat'";        Val=0}
    @{Row=20; Name="''Parabacteroides_johnsonii_DSM_18315.mat'";         Val=0.813}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.Name
    $ws.Cells.Item($r, 3).Value = $item.Val
}
